$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the ReadMe note that lived in D1
$ws.Range("D1").ClearContents()

# Prepare the new rows 11 and 12 by copying the formatting from row 10
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("B10").Copy($ws.Range("B11"))
$ws.Range("A10").Copy($ws.Range("A12"))
$ws.Range("B10").Copy($ws.Range("B12"))

# New SCZ donor for row 11
$ws.Range("B11").Value2 = "Br6032"

# Fill in the previously-empty A10 (new NTC donor)
$ws.Range("A10").Value2 = "Br5436"

# New NTC donor for row 11
$ws.Range("A11").Value2 = "Br5931"

# New row 12: NTC donor + repeated SCZ donor
$ws.Range("A12").Value2 = "Br6389"
$ws.Range("B12").Value2 = "Br5746"

# Update the active selection to D1 (matches the saved view state)
$ws.Range("D1").Select()
